# Added a "clear" button workflow equivalent: set the Box quantities that
# were collected for these titles. Row 2 (first title) got a bulk count of
# 200, and a handful of other rows were bumped to 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 200

$rows = @(11, 12, 16, 18, 27, 29, 31, 32, 33, 34, 38, 40, 42, 46, 54, 56, 57, 61, 62, 63, 65, 67, 70, 89, 91, 92, 98, 106, 107, 110, 111)

foreach ($r in $rows) {
    $ws.Range("B$r").Value = 2
}
